# Apply the "Roles, control duplicados y limpieza de credenciales" update.
# Duplicate-control / role cleanup changed the underlying counts, which
# ripples through the three summary sheets: daily global, per-worker and
# per-work-center (CT).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Resumen_diario_global -----------------------------------
# Mesas_registradas, Par_OK and PPI_OK increase for 2025-11 (serial 45977).
$wsGlobal = $wb.Worksheets.Item("Resumen_diario_global")
$wsGlobal.Cells.Item(2, 2).Value = 18   # B2 Mesas_registradas: 15 -> 18
$wsGlobal.Cells.Item(2, 3).Value = 15   # C2 Par_OK: 12 -> 15
$wsGlobal.Cells.Item(2, 5).Value = 15   # E2 PPI_OK: 12 -> 15

# --- Sheet 2: Resumen_por_trabajador -----------------------------------
# Worker 1's Mesas_registradas increases from 5 to 8.
$wsWorker = $wb.Worksheets.Item("Resumen_por_trabajador")
$wsWorker.Cells.Item(2, 3).Value = 8    # C2: 5 -> 8

# --- Sheet 3: Resumen_por_CT -------------------------------------------
# CT 1's Mesas_registradas increases from 5 to 7, and a previously
# duplicated/missing CT "15" now appears, inserted in sorted order
# between CT 10 and CT 24 -- shifting the rows below it down by one.
$wsCT = $wb.Worksheets.Item("Resumen_por_CT")

$wsCT.Cells.Item(2, 3).Value = 7        # C2: 5 -> 7

# Insert a new row before the current row 9 (which holds CT 100) so the
# existing rows 7-9 shift down to 8-10, then fix up the CT numbers.
$wsCT.Rows.Item(9).Insert()

$wsCT.Cells.Item(7, 1).Value = 15       # A7: 24 -> 15
$wsCT.Cells.Item(8, 1).Value = 24       # A8: 50 -> 24

$wsCT.Cells.Item(9, 1).Value = 50       # A9 (new row): CT 50
$wsCT.Cells.Item(9, 2).Value = 45977    # B9 (new row): Fecha
$wsCT.Cells.Item(9, 2).NumberFormat = $wsCT.Cells.Item(8, 2).NumberFormat
$wsCT.Cells.Item(9, 3).Value = 1        # C9 (new row): Mesas_registradas

# Row 10 already holds A10=100, B10=45977, C10=1 after the insert shifted
# the former row 9 down -- nothing further to change there.
